$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture values from the existing row 61 before we touch anything ---
# (these get reused / moved down into the new row 62)
$codeVal   = $ws.Range("G61").Value2      # transaction code, e.g. " XRP/USDT0000008"
$currVal   = $ws.Range("C61").Value2      # "Currency" text, e.g. "        XRP"
$sellVal   = $ws.Range("B60").Value2      # rich "Sell" text (red), same as B62 needs
$profitVal = $ws.Range("K61").Value2      # "Profit(%)" placeholder text
$oldStatus = $ws.Range("H61").Value2      # "IN PROGRESS" -> moves down to row 62

# --- Insert a new row below row 61 (becomes row 62), inheriting formatting ---
$ws.Rows("62:62").Insert()
$ws.Rows("62:62").RowHeight = 14.25

# --- Update row 61: the sell order is now marked DONE, with finalize date and fee ---
# (set first so the new shared-string entries land in the same index order as the
# original authoring session: J61, then F62, E62, D62)
$ws.Range("H61").Value = "DONE"
$ws.Range("I61").Value = 42874.355543981481
$ws.Range("J61").Value = "0.18150000 XRP (0.15%)"

# --- Fill in the new row 62 (the former "in progress" sell order that is now split) ---
$ws.Range("A62").Value = 42874.355543981481
$ws.Range("B62").Value = $sellVal
$ws.Range("C62").Value = $currVal
$ws.Range("F62").Value = "         120 XRP"
$ws.Range("E62").Value = "         0.375  USDT"

# D62's text is purely numeric after trimming ("0.337" with a trailing line break),
# so a direct .Value assignment would get auto-coerced to a number. Build the text
# in a scratch cell (forced to Text format, then reset to Normal so no stray
# number-format survives), then copy just the resulting value/type into D62 so its
# own (wrap-text) formatting, inherited from the row insert, is left untouched.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "              0.337`n"
$scratch.Style = "Normal"
$scratch.Copy()
$d62 = $ws.Range("D62")
$d62.PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

$ws.Range("G62").Value = $codeVal
$ws.Range("H62").Value = $oldStatus
$ws.Range("K62").Value = $profitVal

# --- Move the active selection down to the newly added row, like the author did ---
$ws.Range("I62").Select() | Out-Null
